$d = $word.ActiveDocument

# --- Part 1: three paragraphs whose <w:numPr> has only <w:numId w:val="0"/> need an
# explicit <w:ilvl w:val="0"/> added (ListOutdent on a paragraph that isn't actually
# indented normalizes the numPr to ilvl=0 while preserving numId=0). ---
$d.Paragraphs.Item(39).Range.ListFormat.ListOutdent()
$d.Paragraphs.Item(40).Range.ListFormat.ListOutdent()
$d.Paragraphs.Item(42).Range.ListFormat.ListOutdent()

# --- Part 2: move the _GoBack bookmark from the end of the "Xinyue Zhang and Jian
# Zhang presented..." paragraph to the start of the following (empty) paragraph. ---
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()
$r = $d.Range(2640, 2642)
$d.Bookmarks.Add("_GoBack", $r)
